# Apply cryptos list update (prices + volume%) per commit:
# "Updated cryptos list on Sat Oct  7 11:23:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.985.89"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.641.14"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'212.98"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'0.525"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'23.59"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "'0.260"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'0.0615"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "'0.0883"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("D12").Value = "1.873.87"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "1.645.31"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.576"
$ws.Range("E14").Value = "  +3.80%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'4.10"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "'65.92"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").Value = "27.982.49"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "'234.91"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'4.37"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "'151.21"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'6.98"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "'0.0484"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").Value = "1.421.90"
$ws.Range("E34").Value = "  -3.41%  "
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").Value = "'0.883"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -4.33%  "
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +6.58%  "
$ws.Range("D44").Value = "'66.62"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "1.782.60"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'87.74"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "'7.61"
$ws.Range("E51").Value = "  -1.07%  "
